$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 445
$ws1.Range("F8").Value = 2018
$ws1.Range("F10").Value = 37
$ws1.Range("F11").Value = 35
$ws1.Range("F12").Value = 1600
$ws1.Range("F13").Value = 1600
$ws1.Range("F14").Value = 1328
$ws1.Range("F20").Value = 474
$ws1.Range("F23").Value = 7057
$ws1.Range("F24").Value = 7057
$ws1.Range("F25").Value = 7648
$ws1.Range("F27").Value = 3
$ws1.Range("F28").Value = 182
$ws1.Range("F30").Value = 80
$ws1.Range("F32").Value = 248
$ws1.Range("F33").Value = 172
$ws1.Range("F38").Value = 1391
$ws1.Range("F39").Value = 19
$ws1.Range("F41").Value = 282
$ws1.Range("F45").Value = 313
$ws1.Range("F46").Value = 223
$ws1.Range("F47").Value = 187
$ws1.Range("F49").Value = 137
$ws1.Range("F50").Value = 140

# Sheet: 本地生活 (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 2579
$ws3.Range("F4").Value = 263
$ws3.Range("F5").Value = 125

# Sheet: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 263
$ws4.Range("F7").Value = 125
$ws4.Range("F11").Value = 445
$ws4.Range("F12").Value = 2018
$ws4.Range("F13").Value = 37
$ws4.Range("F14").Value = 35
$ws4.Range("F15").Value = 1600
$ws4.Range("F16").Value = 1600
$ws4.Range("F17").Value = 1328
$ws4.Range("F20").Value = 474
$ws4.Range("F24").Value = 7057
$ws4.Range("F25").Value = 7057
$ws4.Range("F26").Value = 7648
$ws4.Range("F28").Value = 3
$ws4.Range("F29").Value = 80
$ws4.Range("F30").Value = 248
$ws4.Range("F34").Value = 1391
$ws4.Range("F35").Value = 19
$ws4.Range("F38").Value = 282
$ws4.Range("F45").Value = 313
$ws4.Range("F46").Value = 223
$ws4.Range("F48").Value = 137
